# TC01_INS_SpecialTopic-CancerMoonshot.xlsx edits
# - StatQuery (C2): query now ends with a trailing semicolon
# - GrantsTab query (B4): gnt.project_end_date -> gnt.grant_end_date
# - PublicationsTab query (B5): pub.title -> pub.publication_title, and two
#   extra CASE branches added for relative_citation_ratio 1.0 / 2.0
# - re-normalize the (duplicate) wrap-text/size-12 formatting on the edited
#   cells, and move the active selection from C3 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.focus_area LIKE ''%Cancer Moonshot%'';'

$ws.Range("B4").Value = 'SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.focus_area LIKE ''%Cancer Moonshot%''
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;'

$ws.Range("B5").Value = 'SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN ''0''
    WHEN pub.relative_citation_ratio = 7.0 THEN ''7''
    WHEN pub.relative_citation_ratio = 1.0 THEN ''1''
    WHEN pub.relative_citation_ratio = 2.0 THEN ''2''
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
    prg.focus_area LIKE ''%Cancer Moonshot%''
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;'

# Re-apply the wrap/size-12 formatting so identical (duplicate) style
# records collapse back down to the shared one, same as Excel does when the
# cell text is retyped/reformatted.
$ws.Range("B2").Font.Size = 12
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Font.Size = 12
$ws.Range("C2").WrapText = $true
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").WrapText = $true
$ws.Range("B4").Font.Size = 12
$ws.Range("B4").WrapText = $true
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# Move the window/selection: scroll back to the top and select C5 (was C3).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("C5").Select()
